# Fruta / hortaliza, semanal
# Insert 3 new weekly records for "Vega Modelo de Temuco" - Naranja (rows 522:524),
# pushing the existing rows 522:544 down to 525:547.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 522 (inherits formatting from row above).
$ws.Rows("522:524").Insert()

# --- New row 522 ---
$ws.Cells.Item(522, 1).Value = 10
$ws.Cells.Item(522, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(522, 3).Value = "La Araucanía"
$ws.Cells.Item(522, 4).Value = 44516
$ws.Cells.Item(522, 5).Value = 9
$ws.Cells.Item(522, 6).Value = "Fruta"
$ws.Cells.Item(522, 7).Value = 100102
$ws.Cells.Item(522, 8).Value = "Cítricos"
$ws.Cells.Item(522, 9).Value = 100102005
$ws.Cells.Item(522, 10).Value = "Naranja"
$ws.Cells.Item(522, 11).Value = "Navel Late"
$ws.Cells.Item(522, 12).Value = "Primera"
$ws.Cells.Item(522, 13).Value = 8
$ws.Cells.Item(522, 14).Value = 180000
$ws.Cells.Item(522, 15).Value = 190000
$ws.Cells.Item(522, 16).Value = 185000
$ws.Cells.Item(522, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(522, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(522, 19).Value = 462
$ws.Cells.Item(522, 20).Value = 400

# --- New row 523 ---
$ws.Cells.Item(523, 1).Value = 10
$ws.Cells.Item(523, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(523, 3).Value = "La Araucanía"
$ws.Cells.Item(523, 4).Value = 44516
$ws.Cells.Item(523, 5).Value = 9
$ws.Cells.Item(523, 6).Value = "Fruta"
$ws.Cells.Item(523, 7).Value = 100102
$ws.Cells.Item(523, 8).Value = "Cítricos"
$ws.Cells.Item(523, 9).Value = 100102005
$ws.Cells.Item(523, 10).Value = "Naranja"
$ws.Cells.Item(523, 11).Value = "Valencia"
$ws.Cells.Item(523, 12).Value = "Primera"
$ws.Cells.Item(523, 13).Value = 225
$ws.Cells.Item(523, 14).Value = 9000
$ws.Cells.Item(523, 15).Value = 10000
$ws.Cells.Item(523, 16).Value = 9556
$ws.Cells.Item(523, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(523, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(523, 19).Value = 637
$ws.Cells.Item(523, 20).Value = 15

# --- New row 524 ---
$ws.Cells.Item(524, 1).Value = 10
$ws.Cells.Item(524, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(524, 3).Value = "La Araucanía"
$ws.Cells.Item(524, 4).Value = 44516
$ws.Cells.Item(524, 5).Value = 9
$ws.Cells.Item(524, 6).Value = "Fruta"
$ws.Cells.Item(524, 7).Value = 100102
$ws.Cells.Item(524, 8).Value = "Cítricos"
$ws.Cells.Item(524, 9).Value = 100102005
$ws.Cells.Item(524, 10).Value = "Naranja"
$ws.Cells.Item(524, 11).Value = "Valencia"
$ws.Cells.Item(524, 12).Value = "Primera"
$ws.Cells.Item(524, 13).Value = 5
$ws.Cells.Item(524, 14).Value = 180000
$ws.Cells.Item(524, 15).Value = 180000
$ws.Cells.Item(524, 16).Value = 180000
$ws.Cells.Item(524, 17).Value = "$/bins (400 kilos)"
$ws.Cells.Item(524, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(524, 19).Value = 450
$ws.Cells.Item(524, 20).Value = 400
